# Doing Updates for Financials
# - Insert a new "most recent period" column before column D on the SRCE sheet
#   (this shifts the existing D:K data right to E:L).
# - Populate the new column D with the latest period's figures.
# - Correct one historical data point (Capital Expenditures) that shifted into
#   column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRCE")

# 1) Insert a new column before D; this shifts old D:K -> E:L, carrying over
#    values AND formatting/styles automatically.
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D has no explicit formatting yet - clone it
#    from column E (which now holds what used to be column D) so number
#    formats/styles line up with the rest of the table.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)  # xlPasteFormats

# 3) Fill the new column D with the latest period's values, row by row.
$newColumnData = @(
    @{Row=7; Value=43465},
    @{Row=8; Value=257300},
    @{Row=9; Value="NA"},
    @{Row=10; Value="NA"},
    @{Row=12; Value="NA"},
    @{Row=13; Value=0},
    @{Row=14; Value=0},
    @{Row=15; Value=-26200},
    @{Row=17; Value=62900},
    @{Row=18; Value=194400},
    @{Row=20; Value=-89400},
    @{Row=21; Value=137900},
    @{Row=22; Value=0},
    @{Row=23; Value=105000},
    @{Row=24; Value=23500},
    @{Row=25; Value=0},
    @{Row=26; Value=81500},
    @{Row=27; Value=81000},
    @{Row=28; Value=0},
    @{Row=29; Value=900},
    @{Row=30; Value=0},
    @{Row=31; Value=0},
    @{Row=32; Value=89400},
    @{Row=33; Value=81900},
    @{Row=34; Value=0},
    @{Row=35; Value=81900},
    @{Row=38; Value=43465},
    @{Row=41; Value=94900},
    @{Row=42; Value=32600},
    @{Row=43; Value=0},
    @{Row=44; Value=0},
    @{Row=45; Value=0},
    @{Row=46; Value=0},
    @{Row=47; Value=23500},
    @{Row=48; Value=186600},
    @{Row=49; Value=84000},
    @{Row=50; Value=0},
    @{Row=51; Value=0},
    @{Row=52; Value=0},
    @{Row=53; Value=0},
    @{Row=54; Value=6293700},
    @{Row=57; Value=0},
    @{Row=58; Value=0},
    @{Row=59; Value=78600},
    @{Row=60; Value=0},
    @{Row=61; Value=83400},
    @{Row=62; Value=0},
    @{Row=63; Value=0},
    @{Row=64; Value=0},
    @{Row=65; Value=0},
    @{Row=66; Value=5531700},
    @{Row=68; Value=0},
    @{Row=69; Value=0},
    @{Row=70; Value=0},
    @{Row=71; Value=0},
    @{Row=72; Value=399000},
    @{Row=73; Value=0},
    @{Row=74; Value=0},
    @{Row=75; Value=0},
    @{Row=76; Value=762100},
    @{Row=77; Value=0},
    @{Row=80; Value=43465},
    @{Row=81; Value=81900},
    @{Row=83; Value=32800},
    @{Row=84; Value=0},
    @{Row=85; Value=0},
    @{Row=86; Value=0},
    @{Row=87; Value=0},
    @{Row=88; Value=0},
    @{Row=89; Value=144400},
    @{Row=91; Value=-24200},
    @{Row=92; Value=0},
    @{Row=93; Value=0},
    @{Row=94; Value=-442600},
    @{Row=96; Value=-25700},
    @{Row=97; Value=0},
    @{Row=98; Value=0},
    @{Row=99; Value=0},
    @{Row=100; Value=319300},
    @{Row=101; Value=0},
    @{Row=102; Value=21000}
)

foreach ($item in $newColumnData) {
    $ws.Range("D$($item.Row)").Value = $item.Value
}

# 4) Data correction: the "Capital Expenditures" figure that shifted from the
#    old column I into column J was restated from -7300 to -9500.
$ws.Range("J91").Value = -9500

Write-Host "Financials update applied."
